$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final player roster table (header stays the same; rows 2-19 get rewritten)
$data = @(
    @("Payton Pritchard", "PG", "Boston Celtics"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Grayson Allen", "PG,SG", "Phoenix Suns"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Deandre Ayton", "C", "Portland Trail Blazers"),
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Robert Williams III", "C", "Portland Trail Blazers"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Tari Eason", "SF,PF", "Houston Rockets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
